# Update the stochastic-generator sample values ("alpha_zero" linear example).
# The workbook stores these numbers as literal text (shared strings) rather
# than numeric cell values (it was produced by XLSX.jl), so each cell is
# written with a leading quote (forces text entry, like typing '-4.65... in
# Excel) and then restyled back to Normal so the quote character itself does
# not appear in the stored text.

function Set-TextValue {
    param($Range, [string]$Text)
    $Range.Value = "'" + $Text
    $Range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# --- Restricciones_del_follower ---------------------------------------
$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")

Set-TextValue $wsFollower.Range("B2") "-4.657691821664619"
Set-TextValue $wsFollower.Range("D2") "0.24011722556595838"
Set-TextValue $wsFollower.Range("F2") "0.8694817325713169"

Set-TextValue $wsFollower.Range("B3") "0.6576918216646188"
Set-TextValue $wsFollower.Range("D3") "0.1083236165390392"
Set-TextValue $wsFollower.Range("E3") "0.8815947997511641"

Set-TextValue $wsFollower.Range("B4") "-4.4622558915346"
Set-TextValue $wsFollower.Range("D4") "0.1102758390135593"
Set-TextValue $wsFollower.Range("E4") "0.5026320821534822"

Set-TextValue $wsFollower.Range("B5") "3.6605518210954013"
Set-TextValue $wsFollower.Range("D5") "0.3168885247170169"
Set-TextValue $wsFollower.Range("E5") "0.6592204216829478"

Set-TextValue $wsFollower.Range("B6") "-9.066098643186933"
Set-TextValue $wsFollower.Range("D6") "0.4167665579899481"
Set-TextValue $wsFollower.Range("E6") "0.977109604068878"
Set-TextValue $wsFollower.Range("F6") "0.7336423209960407"

# --- Punto_modificado ---------------------------------------------------
$wsPoint = $wb.Worksheets.Item("Punto_modificado")

Set-TextValue $wsPoint.Range("A2") "5.875357499928848"
Set-TextValue $wsPoint.Range("B2") "4.657691821664619"
Set-TextValue $wsPoint.Range("C2") "2.112315956957238"

# --- Vector_bf --------------------------------------------------------
# NOTE: worksheet name lookup is case-insensitive, and this workbook has
# both "Vector_bf" (sheet 5) and "Vector_BF" (sheet 6), which differ only
# by case. Use the sheet index to address each one unambiguously.
$wsBf = $wb.Worksheets.Item(5)   # "Vector_bf"

Set-TextValue $wsBf.Range("A2") "1.0970590052220022"

# --- Vector_BF ----------------------------------------------------------
$wsBF = $wb.Worksheets.Item(6)   # "Vector_BF"

Set-TextValue $wsBF.Range("A2") "-1.3142800010188616"
Set-TextValue $wsBF.Range("A3") "0.9002435759362326"
